$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 66,2
$arr[0,0] = -0.2250723224242941
$arr[0,1] = -0.8235414954970641
$arr[1,0] = -0.4387644391624299
$arr[1,1] = -1.147074961362611
$arr[2,0] = -0.5009110175439728
$arr[2,1] = -1.485441080562453
$arr[3,0] = -0.1531302537635508
$arr[3,1] = -0.573688788488365
$arr[4,0] = -0.1656575655538209
$arr[4,1] = -0.07402481459760665
$arr[5,0] = -0.2963654443817192
$arr[5,1] = -0.9393552922553143
$arr[6,0] = -0.5277503733184525
$arr[6,1] = -1.23504387918215
$arr[7,0] = -0.562776351763446
$arr[7,1] = -1.502304395717307
$arr[8,0] = -0.4967388911302477
$arr[8,1] = -0.6338882196675701
$arr[9,0] = -0.2914814142137754
$arr[9,1] = -1.069564228628162
$arr[10,0] = -0.2314283557367993
$arr[10,1] = -0.5592037072918058
$arr[11,0] = -0.04778473481036863
$arr[11,1] = -0.1046049227360559
$arr[12,0] = -0.2938613274273816
$arr[12,1] = -0.4221983671803622
$arr[13,0] = -0.1429591273362837
$arr[13,1] = 0.000642066226674419
$arr[14,0] = -0.1893294082381123
$arr[14,1] = -0.4684428268900849
$arr[15,0] = 0.05955896605207452
$arr[15,1] = 0.50688396453999
$arr[16,0] = 0.03705577282857093
$arr[16,1] = 0.6247782069348545
$arr[17,0] = 0.06061489828042262
$arr[17,1] = 0.5618666914615709
$arr[18,0] = -0.1571418069833503
$arr[18,1] = -0.3964320809537481
$arr[19,0] = -0.003471962389699979
$arr[19,1] = 0.1883330593654449
$arr[20,0] = 0.03167019222035432
$arr[20,1] = 0.4756834250290031
$arr[21,0] = -0.006387495855463117
$arr[21,1] = 0.3088153795977731
$arr[22,0] = 0.8133986947206518
$arr[22,1] = 2.032610192429753
$arr[23,0] = 0.125070091254892
$arr[23,1] = 0.7297210752058592
$arr[24,0] = 0.1411341267904336
$arr[24,1] = 0.6497430028463023
$arr[25,0] = 0.09689679049832868
$arr[25,1] = 0.7293672251835698
$arr[26,0] = 0.3134564689051969
$arr[26,1] = 1.233745672600365
$arr[27,0] = 0.6786060401717904
$arr[27,1] = 2.156780179750012
$arr[28,0] = 0.2093136870834215
$arr[28,1] = 0.8412898168301894
$arr[29,0] = 0.01529305187396975
$arr[29,1] = 0.5569674790314447
$arr[30,0] = 0.1811026071771276
$arr[30,1] = 1.057994511581706
$arr[31,0] = 0.1182055391581174
$arr[31,1] = 0.8624821284206031
$arr[32,0] = 0.08768050121899851
$arr[32,1] = 0.4927428815764155
$arr[33,0] = 0.4353820777895407
$arr[33,1] = 1.411243237406499
$arr[34,0] = 0.2912350789677989
$arr[34,1] = 0.9632579243908342
$arr[35,0] = 0.05042502672056587
$arr[35,1] = 0.1470352512500037
$arr[36,0] = 0.3450099733449717
$arr[36,1] = 1.719393195827562
$arr[37,0] = -0.08945859112374727
$arr[37,1] = -0.2339953357181078
$arr[38,0] = 0.1644749660519773
$arr[38,1] = 0.786559508789839
$arr[39,0] = -0.1658191690439932
$arr[39,1] = 0.6223379955691422
$arr[40,0] = 0.2550508305455233
$arr[40,1] = 1.558656059231163
$arr[41,0] = 0.1671319211945836
$arr[41,1] = 0.9316666398835916
$arr[42,0] = -0.1038040671467589
$arr[42,1] = 0.2299779245084234
$arr[43,0] = -0.1047933675262886
$arr[43,1] = 0.1371378966828109
$arr[44,0] = -0.1925832340903354
$arr[44,1] = -0.5585989308011797
$arr[45,0] = -0.1900733665003345
$arr[45,1] = -0.4792173105148722
$arr[46,0] = -0.2281223999680888
$arr[46,1] = -0.5954019046430893
$arr[47,0] = -0.2203032891604934
$arr[47,1] = -0.6384514050277401
$arr[48,0] = -0.1494118902292621
$arr[48,1] = -0.2895409514732279
$arr[49,0] = -0.2130496504420102
$arr[49,1] = -0.6766447720347262
$arr[50,0] = -0.2130496504420102
$arr[50,1] = -0.6766447720347262
$arr[51,0] = -0.2040384834060711
$arr[51,1] = -0.4508419877414442
$arr[52,0] = -0.2095110159062245
$arr[52,1] = -0.5842841537415164
$arr[53,0] = -0.1720812250748444
$arr[53,1] = -0.4355638690044962
$arr[54,0] = -0.158755944584256
$arr[54,1] = -0.481504745948469
$arr[55,0] = -0.1972685240135255
$arr[55,1] = -0.4368972990561003
$arr[56,0] = -0.1813712751261744
$arr[56,1] = -0.3844918846135882
$arr[57,0] = -0.2222104964209077
$arr[57,1] = -0.6634528545323003
$arr[58,0] = -0.2095291394061949
$arr[58,1] = -0.6312153870798277
$arr[59,0] = -0.2722058757751286
$arr[59,1] = -0.35158947413762
$arr[60,0] = -0.1606004435188407
$arr[60,1] = -0.09261315212649859
$arr[61,0] = -0.3529602621162478
$arr[61,1] = -1.121014838932152
$arr[62,0] = -0.2577444427052165
$arr[62,1] = -0.5752547464388464
$arr[63,0] = -0.1998220054726317
$arr[63,1] = -0.5473207601985652
$arr[64,0] = -0.1203783885389504
$arr[64,1] = -0.1089765748813438
$arr[65,0] = -0.06219604907842365
$arr[65,1] = -0.06938642046030193

$ws.Range("A2:B67").Value = $arr

$wb.Save()